$d = $word.ActiveDocument

$replacements = @(
    @("81×37=2997", "63×15=945"),
    @("67×11=737", "61×61=3721"),
    @("53×97=5141", "60×75=4500"),
    @("40×26=1040", "23×93=2139"),
    @("70×79=5530", "62×88=5456"),
    @("99×95=9405", "49×99=4851"),
    @("11×97=1067", "26×12=312"),
    @("19×89=1691", "99×73=7227"),
    @("98×51=4998", "63×33=2079"),
    @("69×23=1587", "34×75=2550"),
    @("72×63=4536", "34×91=3094"),
    @("96×38=3648", "20×18=360"),
    @("77×92=7084", "69×54=3726"),
    @("63×44=2772", "91×77=7007"),
    @("88×26=2288", "37×38=1406"),
    @("41×20=820", "24×31=744"),
    @("18×88=1584", "25×81=2025"),
    @("46×32=1472", "54×72=3888"),
    @("93×79=7347", "34×76=2584"),
    @("78×51=3978", "66×55=3630"),
    @("96×47=4512", "41×71=2911"),
    @("56×73=4088", "66×11=726"),
    @("44×42=1848", "35×35=1225"),
    @("50×93=4650", "61×62=3782"),
    @("87×36=3132", "36×13=468")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Write-Host "Done applying replacements"
